$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing Excel row 214 (pushes 214..314 down
# to 215..315, growing the used range from A1:R314 to A1:R315).
$ws.Rows(214).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(214, 1).Value = 3
$ws.Cells.Item(214, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(214, 3).Value = "Coquimbo"
$ws.Cells.Item(214, 4).Value = 44523
$ws.Cells.Item(214, 5).Value = 5
$ws.Cells.Item(214, 6).Value = 100112021
$ws.Cells.Item(214, 7).Value = "Ají"
$ws.Cells.Item(214, 8).Value = "Americana (o)"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 35
$ws.Cells.Item(214, 11).Value = 31000
$ws.Cells.Item(214, 12).Value = 31000
$ws.Cells.Item(214, 13).Value = 31000
$ws.Cells.Item(214, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(214, 15).Value = "Limache"
$ws.Cells.Item(214, 16).Value = 2067
$ws.Cells.Item(214, 17).Value = 15
$ws.Cells.Item(214, 18).Value = "Hortaliza"
